$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 228.91667
$ws.Range("I12").Value = 127.55556
$ws.Range("K12").Value = 127.55556
$ws.Range("M12").Value = 42.44444
$ws.Range("H40").Value = 2536.182
$ws.Range("I40").Value = 2199
$ws.Range("J40").Value = 2611.111
$ws.Range("K40").Value = 2199
$ws.Range("L40").Value = 2611.111
$ws.Range("M40").Value = -2024
$ws.Range("N40").Value = -2961.111
$ws.Range("H58").Value = 3491.1428
$ws.Range("I58").Value = 284.5
$ws.Range("J58").Value = 7766.6665
$ws.Range("K58").Value = 853.5
$ws.Range("L58").Value = 23299.9995
$ws.Range("M58").Value = -703.5
$ws.Range("N58").Value = -23599.9995
$ws.Range("H70").Value = 851581.0600000001
$ws.Range("I70").Value = 2551768.5
$ws.Range("J70").Value = 1487.375
$ws.Range("K70").Value = 7655305.5
$ws.Range("L70").Value = 4462.125
$ws.Range("M70").Value = -7655035.5
$ws.Range("N70").Value = -5002.125
$ws.Range("H73").Value = 851581.0600000001
$ws.Range("I73").Value = 2551768.5
$ws.Range("J73").Value = 1487.375
$ws.Range("K73").Value = 7655305.5
$ws.Range("L73").Value = 4462.125
$ws.Range("M73").Value = -7654369.5
$ws.Range("N73").Value = -6334.125
$ws.Range("H101").Value = 378.16666
$ws.Range("I101").Value = 378.16666
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 1134.49998
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 487.5000199999999
$ws.Range("N101").ClearContents()
$ws.Range("H106").Value = 59494.945
$ws.Range("I106").Value = 3422.0715
$ws.Range("K106").Value = 3422.0715
$ws.Range("M106").Value = -2791.0715
$ws.Range("H137").Value = 1497.8334
$ws.Range("I137").Value = 1465.9259
$ws.Range("K137").Value = 4397.7777
$ws.Range("M137").Value = -1847.7777

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 67000
$ws.Range("J51").Value = 67000
$ws.Range("L51").Value = 67000
$ws.Range("N51").Value = -67982
$ws.Range("H99").Value = 1033.4736
$ws.Range("I99").Value = 1000.9375
$ws.Range("K99").Value = 1000.9375
$ws.Range("M99").Value = 497.0625
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H131").Value = 73499.5
$ws.Range("J131").Value = 73499.5
$ws.Range("L131").Value = 73499.5
$ws.Range("N131").Value = -83579.5
$ws.Range("H134").Value = 1161.7368
$ws.Range("I134").Value = 1059.6111
$ws.Range("K134").Value = 3178.8333
$ws.Range("M134").Value = -643.8333000000002

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H98").Value = 80944.5
$ws.Range("I98").Value = 77500
$ws.Range("J98").Value = 84389
$ws.Range("K98").Value = 77500
$ws.Range("L98").Value = 84389
$ws.Range("M98").Value = -75254
$ws.Range("N98").Value = -88881
$ws.Range("H99").Value = 3471.1667
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 3471.1667
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 3471.1667
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -6467.1667
$ws.Range("H105").Value = 1745
$ws.Range("I105").Value = 1745
$ws.Range("K105").Value = 1745
$ws.Range("M105").Value = 2
$ws.Range("H107").Value = 1677.7941
$ws.Range("I107").Value = 1601.72
$ws.Range("J107").Value = 1889.1111
$ws.Range("K107").Value = 1601.72
$ws.Range("L107").Value = 1889.1111
$ws.Range("M107").Value = 318.28
$ws.Range("N107").Value = -5729.1111
$ws.Range("H126").Value = 3471.1667
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 3471.1667
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 10413.5001
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -15353.5001
$ws.Range("H134").Value = 1581.3
$ws.Range("I134").Value = 1223.7222
$ws.Range("K134").Value = 3671.1666
$ws.Range("M134").Value = -1136.1666

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 40000
$ws.Range("J37").Value = 40000
$ws.Range("L37").Value = 120000
$ws.Range("N37").Value = -120224
$ws.Range("H92").Value = 103.8
$ws.Range("J92").Value = 94.5
$ws.Range("L92").Value = 283.5
$ws.Range("N92").Value = -2779.5
$ws.Range("H98").Value = 4229.2
$ws.Range("I98").Value = 11522.8
$ws.Range("J98").Value = 582.4
$ws.Range("K98").Value = 34568.39999999999
$ws.Range("L98").Value = 1747.2
$ws.Range("M98").Value = -33070.39999999999
$ws.Range("N98").Value = -4743.2

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 57529.5
$ws.Range("J15").Value = 57529.5
$ws.Range("L15").Value = 57529.5
$ws.Range("N15").Value = -58105.5
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H80").Value = 2351.3076
$ws.Range("I80").Value = 2668.4546
$ws.Range("K80").Value = 2668.4546
$ws.Range("M80").Value = -1670.4546
$ws.Range("H81").Value = 57529.5
$ws.Range("J81").Value = 57529.5
$ws.Range("L81").Value = 57529.5
$ws.Range("N81").Value = -59525.5
$ws.Range("H83").Value = 2351.3076
$ws.Range("I83").Value = 2668.4546
$ws.Range("K83").Value = 13342.273
$ws.Range("M83").Value = -8350.273000000001
$ws.Range("H84").Value = 57529.5
$ws.Range("J84").Value = 57529.5
$ws.Range("L84").Value = 172588.5
$ws.Range("N84").Value = -182572.5
$ws.Range("H97").Value = 918.5714
$ws.Range("I97").Value = 903.1667
$ws.Range("J97").Value = 1011
$ws.Range("K97").Value = 903.1667
$ws.Range("L97").Value = 1011
$ws.Range("M97").Value = -407.1667
$ws.Range("N97").Value = -2003
$ws.Range("H102").Value = 1746.6666
$ws.Range("I102").Value = 1441.7333
$ws.Range("J102").Value = 2509
$ws.Range("K102").Value = 1441.7333
$ws.Range("L102").Value = 2509
$ws.Range("M102").Value = 180.2666999999999
$ws.Range("N102").Value = -5753
$ws.Range("H113").Value = 1430.6666
$ws.Range("I113").Value = 1456.8
$ws.Range("J113").Value = 1300
$ws.Range("K113").Value = 1456.8
$ws.Range("L113").Value = 1300
$ws.Range("M113").Value = 713.2
$ws.Range("N113").Value = -5640
$ws.Range("H132").Value = 4150.107
$ws.Range("I132").Value = 3610.35
$ws.Range("K132").Value = 10831.05
$ws.Range("M132").Value = -8301.049999999999

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 642.5333000000001
$ws.Range("I22").Value = 416.75
$ws.Range("K22").Value = 416.75
$ws.Range("M22").Value = -121.75
$ws.Range("H27").Value = 642.5333000000001
$ws.Range("I27").Value = 416.75
$ws.Range("K27").Value = 416.75
$ws.Range("M27").Value = -309.75
$ws.Range("H46").Value = 1507
$ws.Range("I46").Value = 851.8570999999999
$ws.Range("J46").Value = 1690.44
$ws.Range("K46").Value = 851.8570999999999
$ws.Range("L46").Value = 1690.44
$ws.Range("M46").Value = -663.8570999999999
$ws.Range("N46").Value = -2066.44
$ws.Range("H57").Value = 30000
$ws.Range("I57").Value = 30000
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 30000
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -29434
$ws.Range("N57").ClearContents()
$ws.Range("H122").Value = 3159
$ws.Range("I122").Value = 2272.75
$ws.Range("K122").Value = 6818.25
$ws.Range("M122").Value = -4368.25
$ws.Range("H132").Value = 4819.2256
$ws.Range("I132").Value = 2694.6191
$ws.Range("J132").Value = 9280.9
$ws.Range("K132").Value = 8083.8573
$ws.Range("L132").Value = 27842.7
$ws.Range("M132").Value = -5553.8573
$ws.Range("N132").Value = -32902.7
$ws.Range("H136").Value = 2053.0908
$ws.Range("J136").Value = 2399.5
$ws.Range("L136").Value = 7198.5
$ws.Range("N136").Value = -12298.5

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 5131.8184
$ws.Range("I100").Value = 6831.375
$ws.Range("J100").Value = 599.6667
$ws.Range("K100").Value = 13662.75
$ws.Range("L100").Value = 1199.3334
$ws.Range("M100").Value = -13121.75
$ws.Range("N100").Value = -2281.3334
$ws.Range("H106").Value = 25000
$ws.Range("J106").Value = 25000
$ws.Range("L106").Value = 25000
$ws.Range("N106").Value = -27524
$ws.Range("H107").Value = 567.4286
$ws.Range("I107").Value = 562.1667
$ws.Range("K107").Value = 1686.5001
$ws.Range("M107").Value = 233.4999
$ws.Range("H109").Value = 75749.75
$ws.Range("J109").Value = 75749.75
$ws.Range("L109").Value = 75749.75
$ws.Range("N109").Value = -78523.75
$ws.Range("H113").Value = 3134.5
$ws.Range("I113").Value = 549.2857
$ws.Range("K113").Value = 1647.8571
$ws.Range("M113").Value = 522.1428999999998
$ws.Range("H132").Value = 7303.846
$ws.Range("I132").Value = 8564.429
$ws.Range("K132").Value = 25693.287
$ws.Range("M132").Value = -23163.287
$ws.Range("H136").Value = 3328.087
$ws.Range("I136").Value = 1297.1666
$ws.Range("J136").Value = 10639.4
$ws.Range("K136").Value = 3891.4998
$ws.Range("L136").Value = 31918.2
$ws.Range("M136").Value = -1341.4998
$ws.Range("N136").Value = -37018.2
